$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("8 Dena Demas")

# Update the "Mentor Level" (G column) values to their revised state.
# Empty string clears the cell back to blank; "L"/"M" writes the
# corresponding shared-string value back into the cell.
$ws.Range("G3").Value  = ""
$ws.Range("G5").Value  = ""
$ws.Range("G6").Value  = ""
$ws.Range("G7").Value  = "L"
$ws.Range("G8").Value  = ""
$ws.Range("G10").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("G14").Value = "M"
$ws.Range("G16").Value = "L"
$ws.Range("G17").Value = "M"
$ws.Range("G21").Value = ""
$ws.Range("G22").Value = "L"
$ws.Range("G24").Value = "M"
$ws.Range("G25").Value = ""
$ws.Range("G26").Value = "M"

# Move the on-screen selection to A10:K10, matching the saved cursor
# position captured in the workbook.
$ws.Activate() | Out-Null
$ws.Range("A10:K10").Select() | Out-Null
